# Add Behav Proc in press
# - Replace the "why" link for the CIVN2020 row (E2) which used to point to
#   the CIVN2020 committee page, with a new DOI link (no hyperlink needed
#   there anymore, the ISEP hyperlink on E3 is kept as-is).
# - Remove the old E2 hyperlink; the cell text itself is overwritten.
# - Leave the current selection on E10 (last used cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink that lives on E2 (the civn2020.com committee link),
# keeping the other hyperlink (E3, ISEP congreso2020 page) untouched.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$E$2') {
        $hl.Delete()
        break
    }
}

# Replace the cell's text with the new DOI reference link text.
$ws.Range("E2").Value = "http://doi.org/10.17605/OSF.IO/5BWNX"

# Match the final selection recorded in the saved workbook.
$ws.Range("E10").Select()
